# Update CO2 price results for the relevant years.
# Sheets are named after the target year and hold a single value in A2.
$wb = $excel.ActiveWorkbook

$updates = @{
    "2025" = 57
    "2030" = 195
    "2040" = 355
    "2045" = 355
    "2050" = 355
}

foreach ($key in $updates.Keys) {
    $sheetName = [string]$key
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A2").Value = $updates[$key]
}
